# Manual VIES search included in error handling for API, readme added
#
# The source workbook is a VAT-number / company-name lookup table used by an
# automated VIES-API lookup script. This edit appends two more rows that were
# resolved manually (VIES API failure fallback), and drops the now-unused
# helper columns C:F (they were scratch columns, never populated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 fresh rows below the existing data (rows 9-16), inheriting the
# formatting (style) of row 8 directly above them - this mirrors the way
# Excel's own "insert row" behaviour carries formatting down, including for
# the soon-to-be-deleted helper columns C:F.
$ws.Rows("9:16").Insert()

# Two new VAT Number / Company Name pairs, resolved manually via the VIES
# website after the API lookup failed for them.
$ws.Range("A9").Value = "154186115"
$ws.Range("B9").Value = "RANDSTAD PUBLIC SERVICES LIMITED"

$ws.Range("A10").Value = "492 4357 26"
$ws.Range("B10").Value = "MANPOWER UK HOLDINGS LIMITED"

# Rows 11-16 stay blank (reserved/padding rows), same as the authored file.

# The helper columns C:F were always empty scratch columns - remove them now
# that they're no longer needed.
$ws.Range("C:F").Delete()

# Match the saved selection/view state of the authored workbook.
$ws.Range("F18").Select()

Write-Output "edit applied"
